$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = @{ D = "24.824.79";    E = "  -0.55%  " }
    3  = @{ D = "1.654.59";     E = "  -1.73%  " }
    4  = @{               E = "  +0.29%  " }
    5  = @{ D = "311.46";       E = "  +0.52%  " }
    6  = @{ D = "1.002";        E = "  +0.55%  " }
    7  = @{ D = "0.3620";       E = "  -1.63%  " }
    8  = @{ D = "47.15";        E = "  -0.93%  " }
    9  = @{ D = "0.3238";       E = "  -4.05%  " }
    10 = @{ D = "1.125";        E = "  -4.70%  " }
    11 = @{ D = "0.07026";      E = "  -4.31%  " }
    12 = @{ D = "1.002";        E = "  +0.49%  " }
    13 = @{ D = "6.014";        E = "  -3.01%  " }
    14 = @{ D = "19.36";        E = "  -6.03%  " }
    15 = @{ D = "1.656.60";     E = "  -1.80%  " }
    16 = @{ D = "6.554";        E = "  -4.44%  " }
    17 = @{ D = "0.00001040";   E = "  -5.88%  " }
    18 = @{ D = "0.06585";      E = "  -0.35%  " }
    19 = @{ D = "1.002";        E = "  +0.62%  " }
    20 = @{ D = "78.59";        E = "  -4.73%  " }
    21 = @{ D = "5.868";        E = "  -5.49%  " }
    22 = @{ D = "15.63";        E = "  -7.63%  " }
    23 = @{ D = "12.47";        E = "  -1.42%  " }
    24 = @{ D = "24.805.77";    E = "  -0.27%  " }
    25 = @{ D = "2.429";        E = "  +0.05%  " }
    26 = @{ D = "2.420";        E = "  -10.92%  " }
    27 = @{ D = "147.22";       E = "  -2.33%  " }
    28 = @{ D = "18.49";        E = "  -7.11%  " }
    29 = @{ D = "1.836.31";     E = "  -2.07%  " }
    30 = @{ D = "124.85";       E = "  -4.41%  " }
    31 = @{ D = "1.190";        E = "  -7.27%  " }
    32 = @{ D = "4.077";        E = "  -1.97%  " }
    33 = @{ D = "5.705";        E = "  -12.85%  " }
    34 = @{ D = "0.08438";      E = "  -2.07%  " }
    35 = @{ D = "1.654";        E = "  -4.57%  " }
    36 = @{ D = "12.15";        E = "  -10.22%  " }
    37 = @{ D = "1.280";        E = "  +2.66%  " }
    38 = @{ D = "5.143";        E = "  -5.82%  " }
    39 = @{ D = "0.02243";      E = "  -4.88%  " }
    40 = @{ D = "0.06015";      E = "  -7.33%  " }
    41 = @{ D = "8.233";        E = "  -6.23%  " }
    42 = @{ D = "0.2056";       E = "  -5.48%  " }
    43 = @{ D = "1.002";        E = "  +0.60%  " }
    44 = @{ D = "0.5885";       E = "  -6.51%  " }
    45 = @{ D = "3.762";        E = "  -0.98%  " }
    46 = @{ D = "12.69";        E = "  -5.99%  " }
    47 = @{ D = "0.5576";       E = "  -7.04%  " }
    48 = @{ D = "124.61";       E = "  -1.12%  " }
    49 = @{ D = "1.930";        E = "  -5.92%  " }
    50 = @{ D = "0.06967";      E = "  -2.89%  " }
    51 = @{ D = "1.189";        E = "  -0.97%  " }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey("D")) {
        # Leading apostrophe forces the numeric-looking price string to stay
        # text (matching the original inlineStr/text cell), instead of being
        # auto-converted to a Number by Excel's value parser.
        $ws.Range("D$row").Value = "'" + $vals["D"]
    }
    $ws.Range("E$row").Value = $vals["E"]
}
